# Adding the changes we made on may 9th
#
# The original sheet had header row 1 ("x","y","z") followed by 20 data
# rows (rows 2-21). This change inserts 5 new data rows at the top
# (pushing the existing 20 rows down to rows 7-26) and appends 5 new
# data rows at the bottom (rows 27-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift the existing data rows (2..21) down by 5 rows -----
# Walk from the bottom up so we never clobber a row before reading it.
for ($r = 21; $r -ge 2; $r--) {
    $dest = $r + 5
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value()
}

# --- Step 2: fill the 5 newly-freed rows (2..6) at the top -----------
$ws.Cells.Item(2, 1).Value = 0.0542142912745475
$ws.Cells.Item(2, 2).Value = -0.0145080499351024
$ws.Cells.Item(2, 3).Value = 0.0216857157647609

$ws.Cells.Item(3, 1).Value = -0.0039706239476799
$ws.Cells.Item(3, 2).Value = -0.0080939643085002
$ws.Cells.Item(3, 3).Value = 0.0329867228865623

$ws.Cells.Item(4, 1).Value = -0.0224492978304624
$ws.Cells.Item(4, 2).Value = 0.0106901414692401
$ws.Cells.Item(4, 3).Value = 0.0164933614432811

$ws.Cells.Item(5, 1).Value = -0.0187841057777404
$ws.Cells.Item(5, 2).Value = -0.0056505035609006
$ws.Cells.Item(5, 3).Value = 0.0255036242306232

$ws.Cells.Item(6, 1).Value = -0.008552113547921099
$ws.Cells.Item(6, 2).Value = -0.001527163083665
$ws.Cells.Item(6, 3).Value = 0.0328340083360672

# --- Step 3: append 5 new data rows (27..31) at the bottom -----------
$ws.Cells.Item(27, 1).Value = 0.0029016099870204
$ws.Cells.Item(27, 2).Value = 0.0401643887162208
$ws.Cells.Item(27, 3).Value = 0.0038179077673703

$ws.Cells.Item(28, 1).Value = 0.0178678091615438
$ws.Cells.Item(28, 2).Value = 0.0429132841527462
$ws.Cells.Item(28, 3).Value = 0.027030786499381

$ws.Cells.Item(29, 1).Value = -0.0041233403608202
$ws.Cells.Item(29, 2).Value = -0.0798706337809562
$ws.Cells.Item(29, 3).Value = -0.0476474873721599

$ws.Cells.Item(30, 1).Value = -0.0029016099870204
$ws.Cells.Item(30, 2).Value = 0.0207694191485643
$ws.Cells.Item(30, 3).Value = 0.0029016099870204

$ws.Cells.Item(31, 1).Value = -0.007177666760981
$ws.Cells.Item(31, 2).Value = -0.0045814891345798
$ws.Cells.Item(31, 3).Value = 0.0218384321779012
